# [FEATURE] Adds more preferences to the tictactoe example
#
# Inserts two new numeric columns (G, H) before the existing "comment" text
# column, which pushes that column from G to I, then fills in the new
# columns for every existing row and appends ten new rows (17-26)
# describing additional preference runs that were tried.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Insert two blank columns at G:H - this shifts the old column G (the text
# comment column) to column I, carrying its custom width along with it.
$ws.Range("G1:H1").EntireColumn.Insert()

# Fill the new G/H columns with 10 for every existing data row (2-16).
for ($r = 2; $r -le 16; $r++) {
    $ws.Range("G$r").Value = 10
    $ws.Range("H$r").Value = 10
}

# New rows 17-19 (written before the new header cells below, so the new
# comment strings they introduce land at the shared-string indices the
# original file uses before the "N in 1. L" / "N in 2. L" header strings).
$ws.Range("A17").Value = 1.8
$ws.Range("B17").Value = 0.3
$ws.Range("C17").Value = 250
$ws.Range("D17").Value = 20
$ws.Range("E17").Value = 0
$ws.Range("G17").Value = 10
$ws.Range("H17").Value = 10
$ws.Range("I17").Value = "nach ca 38 M Comp. => 1.0 rating"

$ws.Range("A18").Value = 1.8
$ws.Range("B18").Value = 0.3
$ws.Range("C18").Value = 250
$ws.Range("D18").Value = 30
$ws.Range("E18").Value = 0
$ws.Range("G18").Value = 10
$ws.Range("H18").Value = 10
$ws.Range("I18").Value = "nach ca 50 M Comp. => 1.0 rating"

$ws.Range("A19").Value = 1.8
$ws.Range("B19").Value = 0.3
$ws.Range("C19").Value = 250
$ws.Range("D19").Value = 15
$ws.Range("E19").Value = 0
$ws.Range("F19").Value = 12
$ws.Range("G19").Value = 10
$ws.Range("H19").Value = 10
$ws.Range("I19").Value = "nach ca 30 M Comp. => 1.0 rating"

# New header cells for the two inserted numeric columns.
$ws.Range("G1").Value = "N in 1. L"
$ws.Range("H1").Value = "N in 2. L"

$ws.Range("A20").Value = 1.8
$ws.Range("B20").Value = 0.3
$ws.Range("C20").Value = 250
$ws.Range("D20").Value = 25
$ws.Range("E20").Value = 10
$ws.Range("G20").Value = 10
$ws.Range("H20").Value = 10
$ws.Range("I20").Value = "nach ca 50 M Comp. => 1.0 rating"

$ws.Range("A21").Value = 1.8
$ws.Range("B21").Value = 0.3
$ws.Range("C21").Value = 250
$ws.Range("D21").Value = 10
$ws.Range("E21").Value = 0
$ws.Range("G21").Value = 10
$ws.Range("H21").Value = 10
$ws.Range("I21").Value = "nach ca 50 M Comp. => 1.0 rating"

$ws.Range("A22").Value = 1.8
$ws.Range("B22").Value = 0.3
$ws.Range("C22").Value = 250
$ws.Range("D22").Value = 25
$ws.Range("E22").Value = 0
$ws.Range("G22").Value = 9
$ws.Range("H22").Value = 9
$ws.Range("I22").Value = "nach ca 38 M Comp. => 1.0 rating"

$ws.Range("A23").Value = 1.8
$ws.Range("B23").Value = 0.3
$ws.Range("C23").Value = 250
$ws.Range("D23").Value = 25
$ws.Range("E23").Value = 0
$ws.Range("F23").Value = 13
$ws.Range("G23").Value = 5
$ws.Range("H23").Value = 5
$ws.Range("I23").Value = "nach ca 25 M Comp. => 1.0 rating (8 min)"

$ws.Range("A24").Value = 1.8
$ws.Range("B24").Value = 0.3
$ws.Range("C24").Value = 250
$ws.Range("D24").Value = 25
$ws.Range("E24").Value = 0
$ws.Range("G24").Value = 5
$ws.Range("H24").Value = 3
$ws.Range("I24").Value = "nach ca 40 M Comp. => 1.0 rating"

$ws.Range("A25").Value = 1.8
$ws.Range("B25").Value = 0.3
$ws.Range("C25").Value = 250
$ws.Range("D25").Value = 25
$ws.Range("E25").Value = 0
$ws.Range("G25").Value = 3
$ws.Range("H25").Value = 3
$ws.Range("I25").Value = "nach ca 75 M Comp. => 1.0 rating"

$ws.Range("A26").Value = 1.8
$ws.Range("B26").Value = 0.3
$ws.Range("C26").Value = 250
$ws.Range("D26").Value = 25
$ws.Range("E26").Value = 0
$ws.Range("G26").Value = 5
$ws.Range("H26").Value = 1
$ws.Range("I26").Value = "nach ca 63 M Comp. => 1.0 rating "

# Column A keeps the number-format style used for the rest of the column.
$ws.Range("A17:A26").NumberFormat = "0.00"

# Column widths: the new H (numeric) and I (wide text, inherited from the
# old G) columns get explicit custom widths; G stays at the sheet default.
# (Values chosen land on the closest width the engine's pixel grid supports
# to the authored 11.6640625 / 59.109375 character widths.)
$ws.Columns.Item(8).ColumnWidth = 10.8
$ws.Columns.Item(9).ColumnWidth = 58.3

# Selection moves to J21 in the edited workbook.
$ws.Range("J21").Select()
